$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "User"
$ws.Range("B1").Value = "Mortal"
$ws.Range("C1").Value = "Fact1"
$ws.Range("D1").Value = "Fact2"
$ws.Range("E1").Value = "Fact3"
$ws.Range("F1").Value = "Fact4"
$ws.Range("G1").Value = "Fact5"

# Data rows: User, Mortal, Fact1..Fact5
$data = @(
    @("praveen", "joanne", "p", "2p", "3p", "4p", "5p"),
    @("joanne",  "nick",   "j", "2j", "3j", "4j", "5j"),
    @("nick",    "casper", "n", "2n", "3n", "4n", "5n"),
    @("casper",  "daniel", "c", "2c", "3c", "4c", "5c"),
    @("daniel",  "praveen","d", "2d", "3d", "4d", "5d")
)

$row = 2
foreach ($rowVals in $data) {
    $col = 1
    foreach ($val in $rowVals) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $row++
}

$ws.Range("G7").Select() | Out-Null
